$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("F1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.99699699699699695
$ws.Range("G2").Value = 1

$ws.Range("E3").Value = 0.99708454810495628
$ws.Range("F3").Value = 0.98064516129032253
$ws.Range("G3").Value = 0.97610921501706482

$ws.Columns("F").Select()
